$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new data row (row 80) for DGS's 2021/09/13 report
$row = 80
$prevRow = $row - 1

# Match the date column's existing display format ("yyyy/mm/dd"), but the
# column actually stores the dates as plain text (shared strings), so the
# new date has to be entered as text rather than being auto-converted into
# a date serial number.
$origDateFormat = $ws.Cells.Item($prevRow, 1).NumberFormat
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2021/09/13"
$dateCell.NumberFormat = $origDateFormat

$ws.Cells.Item($row, 2).Value = 208.3
$ws.Cells.Item($row, 3).Value = 214
$ws.Cells.Item($row, 4).Value = 0.85
$ws.Cells.Item($row, 5).Value = 0.84

# Keep the active cell / selection consistent with the appended row
$ws.Range("A81").Select()
